$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "Here is the relevant data on {{ formal_owner }}'s property:"
#    gains a left indent of 720 twips (0.5") -- w:ind w:hanging="0" ->
#    w:ind w:left="720" w:hanging="0"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Here is the relevant data on {{ formal_owner }}", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $para1 = $rng1.Paragraphs(1)
    $para1.LeftIndent = 36
}

# ---------------------------------------------------------------------------
# 2) Typo fix: "...relevant data on most similar property..." ->
#    "...relevant data on the most similar property..." (insert "the ")
# ---------------------------------------------------------------------------
$old2 = "Here is the relevant data on most similar property sold"
$new2 = "Here is the relevant data on the most similar property sold"
$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Typo fix: "...properties by assessed..." -> "...properties be assessed..."
# ---------------------------------------------------------------------------
$old3 = "This assessment is inaccurate. The Michigan Constitution requires that properties by assessed at no more than 50% of their market value."
$new3 = "This assessment is inaccurate. The Michigan Constitution requires that properties be assessed at no more than 50% of their market value."
$rng3 = $d.Content
$rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Collapse the many small runs making up the "These comparable properties..."
#    paragraph (same visible text) into a single run -- done by replacing the
#    stable, already-correct text with itself, which re-serializes the
#    paragraph's runs.
# ---------------------------------------------------------------------------
$text4 = "These comparable properties were selected from all the residential property sales marked arms-length by the Detroit Assessment Division between April 1, 2021 to March 31, 2023."
$rng4 = $d.Content
$rng4.Find.Execute($text4, $true, $false, $false, $false, $false, $true, 1, $false, $text4, 2) | Out-Null
